$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: convert Meta/Meta.AC/Venda/Venda.AC/Sobras/P to real numbers
$ws.Range("B2").Value = 5000
$ws.Range("C2").Value = 5000
$ws.Range("D2").Value = 5000
$ws.Range("E2").Value = 5000
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 100

# Row 3: new record inserted into the list (kept as text, matching source data)
$ws.Range("A3:G3").NumberFormat = "@"
$ws.Range("A3").Value = "03/08/2023"
$ws.Range("B3").Value = "4000.00"
$ws.Range("C3").Value = "9000.00"
$ws.Range("D3").Value = "4000.00"
$ws.Range("E3").Value = "9000.00"
$ws.Range("F3").Value = "0.00"
$ws.Range("G3").Value = "100.00"
$ws.Range("A3:G3").ClearFormats()
